# Apply the data + formatting changes described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C used to hold a repeated shared string ("AlienAnimBP" blueprint path);
# it now holds a simple numeric row index (0-4) instead.
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 4

# Column F (StartingAnimation) values changed for rows 2-5.
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 11
$ws.Range("F6").Value = 0

# Column C is much narrower now that it just stores small numbers.
$ws.Columns("C").ColumnWidth = 12.83

# Update the view: no more horizontal scroll anchored at column C,
# and the remembered selection moved from F10 to C12.
$ws.Range("C12").Select()
